$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Daten Siliziumtombak" block (rows 18-25, columns G-I) ---

# Row 18 is a blank spacer row with custom height
$ws.Rows.Item(18).RowHeight = 15

# Header row 19
$ws.Range("G19").Value = "Daten Siliziumtombak:"
$ws.Range("H19").Value = $null
$ws.Range("I19").Value = $null

# Row 20 (spacer row inside the box, no text)
$ws.Range("G20").Value = $null
$ws.Range("H20").Value = $null
$ws.Range("I20").Value = $null

# Column G labels first
$ws.Range("G21").Value = "Temp"
$ws.Range("G22").Value = "Dichte"
$ws.Range("G23").Value = "kinematische Visko."
$ws.Range("G24").Value = "Dynamische Visko."
$ws.Range("G25").Value = "Oberflächenspannung"

# Column H values
$ws.Range("H21").Value = 920
$ws.Range("H22").Value = 1000
$ws.Range("H23").Formula = "=H24/H22"
$ws.Range("H24").Value = 0.000001
$ws.Range("H25").Value = 0.07

# Column I units
$ws.Range("I21").Value = "[°C]"
$ws.Range("I22").Value = "[kg/m3]"
$ws.Range("I23").Value = "[m2/s]"
$ws.Range("I24").Value = "[kg/m/s] (= 2 mPa*s)"
$ws.Range("I25").Value = "[N/m]"

$ws.Rows.Item(25).RowHeight = 15

# --- Borders around the little boxed table (medium/box style) ---
# Top edge
$ws.Range("G19:I19").Borders.Item(8).Weight = -4138
# Left edge
$ws.Range("G19:G25").Borders.Item(7).Weight = -4138
# Right edge
$ws.Range("I19:I25").Borders.Item(10).Weight = -4138

# Number format for the tiny viscosity values (scientific notation)
$ws.Range("H23:H24").NumberFormat = "0.00E+00"

# Bottom edge
$ws.Range("G25:I25").Borders.Item(9).Weight = -4138
